# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (fund-holder detail) right before the
#   "总计" (totals) summary sheet.
# - Insert a new top row into "总计" for the 2022-Q1 period, pushing the
#   existing 2021-Q4 row down to row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (bold / centered / boxed, matching the other sheets' header
# look).
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (A) mirrors the same bold/centered/boxed style as the
# other sheets' leading "row number" column.
$indexRange = $q1.Range("A2:A4")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Data rows - codes/names/ratios are kept as text (leading "'" preserves
# the literal formatting, e.g. "0.20" would otherwise be coerced to 0.2).
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'160921"
$q1.Range("C2").Value = "大成多策略混合(LOF)"
$q1.Range("D2").Value = "'1.13"
$q1.Range("E2").Value = "'79.19"
$q1.Range("F2").Value = "'4.03"
$q1.Range("G2").Value = "'0.0455"
$q1.Range("H2").Value = 6

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'010690"
$q1.Range("C3").Value = "万家互联互通核心资产量化策略混合A"
$q1.Range("D3").Value = "'0.85"
$q1.Range("E3").Value = "'94.05"
$q1.Range("F3").Value = "'3.07"
$q1.Range("G3").Value = "'0.0261"
$q1.Range("H3").Value = 7

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'010691"
$q1.Range("C4").Value = "万家互联互通核心资产量化策略混合C"
$q1.Range("D4").Value = "'0.20"
$q1.Range("E4").Value = "'94.05"
$q1.Range("F4").Value = "'3.07"
$q1.Range("G4").Value = "'0.0061"
$q1.Range("H4").Value = 7

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计", keeping 2021-Q4 below it.
#    Re-fetch the sheet by name: inserting a sheet before it shifted its
#    tab position, so the earlier $totalSheet handle now resolves to the
#    new "2022-Q1" tab instead.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing 2021-Q4 row from row 2 down to row 3.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.02

# Write the new 2022-Q1 row into row 2.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.08

Write-Output "2022-Q1 sheet added; total sheet updated"
